$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Insert 3 new rows at row 63 (shifts existing rows 63-70 down to 66-73)
$ws.Rows.Item(63).Resize(3).Insert()

# The rows that shifted down (old 65-68, now 68-71) should not carry the
# leftover "E" marker cell from their old position.
$ws.Range("E68:E71").ClearContents()

# New individuals of class Deliver_Stocked_Product, linked to the
# "Spedizione ordine" activity (Activity_1xpl0p1) via has_domainLink.
# Fill column C in reverse row order (65,64,63) so the new shared-string
# entries are minted in the same order as the target file
# (157=EO8199464BBC, 158=EO8199482A97, 159=EO81994532AB).
$ws.Range("A65").Value = "Activity_1xpl0p1"
$ws.Range("B65").Value = "has_domainLink"
$ws.Range("C65").Value = "EO8199464BBC"
$ws.Range("D65").Formula = "=VLOOKUP(A65,Foglio1!A:B,2)"
$ws.Range("E65").Value = ""

$ws.Range("A64").Value = "Activity_1xpl0p1"
$ws.Range("B64").Value = "has_domainLink"
$ws.Range("C64").Value = "EO8199482A97"
$ws.Range("D64").Formula = "=VLOOKUP(A64,Foglio1!A:B,2)"
$ws.Range("E64").Value = ""

$ws.Range("A63").Value = "Activity_1xpl0p1"
$ws.Range("B63").Value = "has_domainLink"
$ws.Range("C63").Value = "EO81994532AB"
$ws.Range("D63").Formula = "=VLOOKUP(A63,Foglio1!A:B,2)"
$ws.Range("E63").Value = ""

# Expand the table (Tabella1) to cover the new rows
$lo = $ws.ListObjects.Item("Tabella1")
$lo.Resize($ws.Range("A1:D71"))

# Recalculate
$wb.Application.Calculate()

# Update sheet view to match final state
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("C63").Select()

# Update workbook window position
$excel.ActiveWindow.Left = 740
$excel.ActiveWindow.Top = 1860
